$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "newasdasd user"

$ws.Range("B2").Select()
